# Implement max item limits per category
$wb = $excel.ActiveWorkbook

# ----- Sheet "items" -----
$items = $wb.Worksheets.Item("items")

# Update min_time (C) / max_time (D) values per row
$items.Range("C2").Value = 2
$items.Range("D2").Value = 5

$items.Range("C3").Value = 2
$items.Range("D3").Value = 5

$items.Range("C4").Value = 2
$items.Range("D4").Value = 5

$items.Range("C5").Value = 2
$items.Range("D5").Value = 5

$items.Range("C6").Value = 2
$items.Range("D6").Value = 5

$items.Range("C7").Value = 2

$items.Range("C12").Value = 5

$items.Range("D13").Value = 10

$items.Range("D14").Value = 10

$items.Range("D17").Value = 10

# Move active selection to D14 on the items sheet
$items.Activate()
$items.Range("D14").Select()

# ----- Sheet "categories" -----
$categories = $wb.Worksheets.Item("categories")

# Remove the old min_items value for the "p1" category (row 2)
$categories.Range("B2").ClearContents()

# Add max_items values for the remaining categories
$categories.Range("C4").Value = 1
$categories.Range("C5").Value = 2

# Remove the old min_items value for the "e1" category (row 6); max_items (C6) stays 1
$categories.Range("B6").ClearContents()

# Match the row heights used by the other data rows on this sheet
$categories.Rows.Item(4).RowHeight = 13.8
$categories.Rows.Item(5).RowHeight = 13.8
$categories.Rows.Item(6).RowHeight = 13.8

# Move active selection to C3 on the categories sheet
$categories.Activate()
$categories.Range("C3").Select()

# Leave the "items" sheet focused/selected as in the original workbook
$items.Activate()
